$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column I: "Deaths"
$ws.Range("I1").Value = "Deaths"
$ws.Range("I1").Style = $ws.Range("A1").Style

# Data for rows 2-10: all zero, with the same style as other data columns
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Value = 0
}
$ws.Range("I2:I10").Style = $ws.Range("B2").Style

# Update selection to I11 (matches the authored sheetView selection)
$ws.Range("I11").Select()
